$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "{'hospital_department', 'Doctor_Name', 'doctor_licence_number'},{'Preventive_care', 'Name_of_clinics', 'List_plan'}"
$ws.Range("D3").Value = "{'Preventive_care', 'Name_of_clinics', 'List_plan'}{'Patient_age', 'Schedule_time', 'Name_of_clinics', 'Customer_phone'},{'Address_Of_clinics', 'Policy_number', 'Preventive_care'}"
$ws.Range("D4").Value = "{'Address_Of_clinics', 'Hospital_Address', 'Acc_type', 'Preventive_care'}"
$ws.Range("D5").Value = "{'hospital_department', 'Doctor_Name', 'doctor_licence_number'}{'hospital_department', 'Doctor_Name', 'doctor_licence_number', 'Hourly_charge_doctor', 'Coverage_policy'},{'Patient_age', 'Hospital_Address', 'Doctor_Name', 'Patient_prior_condition', 'Customer_phone'}"
$ws.Range("D6").Value = "{'Preventive_care', 'Name_of_clinics', 'List_plan'}"
$ws.Range("D7").Value = "{'X,Y_Coordinates'}{'Patient_age', 'Hospital_Address', 'Doctor_Name', 'Patient_prior_condition', 'Customer_phone'}"
$ws.Range("D8").Value = "{'Discharge_amount', 'Hospital_Address', 'Doctor_Name', 'Acc_type', 'Schedule_time'}"
$ws.Range("D9").Value = "{'Doctor_available_time', 'Schedule_time'}"
$ws.Range("D10").Value = "{'Doctor_available_time', 'Schedule_time'}{'hospital_department', 'Doctor_Name', 'doctor_licence_number'},{'hospital_department', 'Doctor_Name', 'doctor_licence_number', 'Hourly_charge_doctor', 'Coverage_policy'}{'Doctor_available_time', 'Schedule_time'}"

$wb.Save()
